# Update the auto date placeholder ("datetimeFigureOut" field) from
# 19-09-2024 to 20-09-2024 across the slide master and every slide layout.
#
# ppPlaceholderDate = 16
$ppPlaceholderDate = 16
$oldDate = "19-09-2024"
$newDate = "20-09-2024"

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout off the master.
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Touching the Guides collection mirrors opening the Guides pane in the
# PowerPoint UI, which is what produced the empty p15:sldGuideLst
# extension on <p:presentation> in the target edit.
try {
    $null = $p.Guides.Count
} catch {
}
